# "Update pawn support panel"
#
# The sheet's last six rows (26-31) held a leftover block of pawn-support
# event options (event ids 301-305, plus the catch-all fallback "999" row).
# The 301-305 rows are being retired, so only the catch-all "999" row
# should remain - it moves up to become the new last row (26) once the
# obsolete rows above it are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 26-30 (ids 301, 302, 303, 304, 305). Row 31 (id 999) shifts
# up to row 26, so the sheet's used range shrinks from A1:D31 to A1:D26
# and the strings that were only referenced by the deleted rows drop out
# of the shared-strings table automatically.
$ws.Range("A26:D30").EntireRow.Delete()

# Leave the selection where the edit finished, on the (now relocated)
# fallback row's effect column.
$ws.Range("F31").Select()
